$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3875
$ws.Range("I82").Value = 1833.3334
$ws.Range("K82").Value = 5500.0002
$ws.Range("M82").Value = -5094.0002

$ws.Range("H85").Value = 3875
$ws.Range("I85").Value = 1833.3334
$ws.Range("K85").Value = 5500.0002
$ws.Range("M85").Value = -4096.0002

$ws.Range("H92").Value = 910.44446
$ws.Range("I92").Value = 1115.8334
$ws.Range("K92").Value = 1115.8334
$ws.Range("M92").Value = 132.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2580.5715
$ws.Range("I2").Value = 1344.3334
$ws.Range("K2").Value = 1344.3334
$ws.Range("M2").Value = -1231.3334

$ws.Range("H5").Value = 507.14285
$ws.Range("I5").Value = 137.5
$ws.Range("K5").Value = 137.5
$ws.Range("M5").Value = -25.5

$ws.Range("H45").Value = 2344.6155
$ws.Range("I45").Value = 1036
$ws.Range("J45").Value = 3162.5
$ws.Range("K45").Value = 1036
$ws.Range("L45").Value = 3162.5
$ws.Range("M45").Value = -659
$ws.Range("N45").Value = -3916.5

$ws.Range("H110").Value = 371.8
$ws.Range("I110").Value = 371.8
$ws.Range("K110").Value = 371.8
$ws.Range("M110").Value = 1673.2

$ws.Range("H116").Value = 2580.5715
$ws.Range("I116").Value = 1344.3334
$ws.Range("K116").Value = 1344.3334
$ws.Range("M116").Value = 949.6666

$ws.Range("H132").Value = 1603
$ws.Range("I132").Value = 1603
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4809
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2279
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2580.5715
$ws.Range("I3").Value = 1344.3334
$ws.Range("K3").Value = 1344.3334
$ws.Range("M3").Value = -1230.3334

$ws.Range("H4").Value = 507.14285
$ws.Range("I4").Value = 137.5
$ws.Range("K4").Value = 137.5
$ws.Range("M4").Value = -22.5

$ws.Range("H20").Value = 998.3333
$ws.Range("I20").Value = 998.3333
$ws.Range("K20").Value = 998.3333
$ws.Range("M20").Value = -751.3333

$ws.Range("H64").Value = 2299.6667
$ws.Range("I64").Value = 999
$ws.Range("J64").Value = 2950
$ws.Range("K64").Value = 999
$ws.Range("L64").Value = 2950
$ws.Range("M64").Value = -774
$ws.Range("N64").Value = -3400

$ws.Range("H67").Value = 2299.6667
$ws.Range("I67").Value = 999
$ws.Range("J67").Value = 2950
$ws.Range("K67").Value = 999
$ws.Range("L67").Value = 2950
$ws.Range("M67").Value = -219
$ws.Range("N67").Value = -4510

$ws.Range("H86").Value = 4005.4285
$ws.Range("I86").Value = 759.75
$ws.Range("J86").Value = 8333
$ws.Range("K86").Value = 759.75
$ws.Range("L86").Value = 8333
$ws.Range("M86").Value = 363.25
$ws.Range("N86").Value = -10579

$ws.Range("H89").Value = 4005.4285
$ws.Range("I89").Value = 759.75
$ws.Range("J89").Value = 8333
$ws.Range("K89").Value = 3798.75
$ws.Range("L89").Value = 41665
$ws.Range("M89").Value = 1817.25
$ws.Range("N89").Value = -52897

$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -80
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 479.625
$ws.Range("I16").Value = 531.5
$ws.Range("J16").Value = 427.75
$ws.Range("K16").Value = 531.5
$ws.Range("L16").Value = 427.75
$ws.Range("M16").Value = -244.5
$ws.Range("N16").Value = -1001.75

$ws.Range("H31").Value = 3707.6667
$ws.Range("I31").Value = 3707.6667
$ws.Range("K31").Value = 3707.6667
$ws.Range("M31").Value = -3412.6667

$ws.Range("H34").Value = 3707.6667
$ws.Range("I34").Value = 3707.6667
$ws.Range("K34").Value = 3707.6667
$ws.Range("M34").Value = -3505.6667

$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("M58").Value = -2797

$ws.Range("H99").Value = 5730
$ws.Range("I99").Value = 4000
$ws.Range("K99").Value = 4000
$ws.Range("M99").Value = -2502

$ws.Range("H113").Value = 479.625
$ws.Range("I113").Value = 531.5
$ws.Range("J113").Value = 427.75
$ws.Range("K113").Value = 531.5
$ws.Range("L113").Value = 427.75
$ws.Range("M113").Value = 1638.5
$ws.Range("N113").Value = -4767.75

$ws.Range("H126").Value = 5730
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("M136").Value = -6450

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 130000
$ws.Range("J141").Value = 130000
$ws.Range("L141").Value = 130000
$ws.Range("N141").Value = -140360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3290.4285
$ws.Range("I46").Value = 390.66666
$ws.Range("J46").Value = 8510
$ws.Range("K46").Value = 1171.99998
$ws.Range("L46").Value = 25530
$ws.Range("M46").Value = -1080.99998
$ws.Range("N46").Value = -25712

$ws.Range("H104").Value = 5623
$ws.Range("I104").Value = 2231
$ws.Range("K104").Value = 6693
$ws.Range("M104").Value = -4072

$ws.Range("H122").Value = 225.5
$ws.Range("I122").Value = 151.5
$ws.Range("J122").Value = 299.5
$ws.Range("K122").Value = 1363.5
$ws.Range("L122").Value = 2695.5
$ws.Range("M122").Value = 1086.5
$ws.Range("N122").Value = -7595.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 357.30768
$ws.Range("I2").Value = 55.25
$ws.Range("J2").Value = 491.55554
$ws.Range("K2").Value = 55.25
$ws.Range("L2").Value = 491.55554
$ws.Range("M2").Value = 57.75
$ws.Range("N2").Value = -717.5555400000001

$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744

$ws.Range("H95").Value = 75000
$ws.Range("J95").Value = 75000
$ws.Range("L95").Value = 75000
$ws.Range("N95").Value = -80492

$ws.Range("H122").Value = 33136.727
$ws.Range("I122").Value = 44975.75
$ws.Range("J122").Value = 1566
$ws.Range("K122").Value = 134927.25
$ws.Range("L122").Value = 4698
$ws.Range("M122").Value = -132477.25
$ws.Range("N122").Value = -9598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2985.7144
$ws.Range("I22").Value = 1966.6666
$ws.Range("J22").Value = 3750
$ws.Range("K22").Value = 1966.6666
$ws.Range("L22").Value = 3750
$ws.Range("M22").Value = -1671.6666
$ws.Range("N22").Value = -4340

$ws.Range("H27").Value = 2985.7144
$ws.Range("I27").Value = 1966.6666
$ws.Range("J27").Value = 3750
$ws.Range("K27").Value = 1966.6666
$ws.Range("L27").Value = 3750
$ws.Range("M27").Value = -1859.6666
$ws.Range("N27").Value = -3964

$ws.Range("H122").Value = 6900
$ws.Range("J122").Value = 7412.5
$ws.Range("L122").Value = 22237.5
$ws.Range("N122").Value = -27137.5
